$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update score values
$ws.Range("B2").Value = 92
$ws.Range("B3").Value = 79
$ws.Range("B4").Value = 59

# Update reason text (replacing the old text, removing leading blank line)
$ws.Range("C2").Value = "The SDE Intern job description is suitable for you because it requires skills in MongoDB, ReactJS, JavaScript, Web Development, and NodeJS, which align with your experience and projects. Additionally, the responsibilities of maintaining code, ensuring scalability, and suggesting new features match well with your project experiences. The high score of 92 indicates a strong match between your skills and the job requirements."

$ws.Range("C3").Value = "The job as a Frontend Engineer Intern (score: 79) is suitable for you because it requires skills such as ReactJS, JavaScript, CSS, Frontend Development which align with your experience in projects like SHAMIYANA APP and SMART SENSING MIDDLEWARE. Your knowledge in these areas will contribute to efficient and visually appealing web design and user experience, making you a valuable addition to the team."

$ws.Range("C4").Value = "The job of NLP engineer is moderately suitable for the candidate. While the candidate has experience in a variety of technologies, including Flutter, Dart, Firebase, NodeJS, ExpressJS, Socket.IO, WebRTC, HTML, CSS, JS, Docker, and ReactJS, they lack specific experience in NLP, Pytorch, Computer Vision, and Python, which are required skills for the job. However, their exposure to various technologies and their ability to learn new ones makes them moderately suitable for the role."

$wb.Save()
